$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 3 new rows before current row 28 (PTPN11), pushing
# PTPN11 -> row 31 and WT1 -> row 32, and leaving rows 28-30 empty for the
# three additional CBL splice-junction rows being added.
$ws.Range("A28:A30").EntireRow.Insert()

# Columns M, N, Q, R, S, T in this sheet are numeric-looking values stored
# as TEXT (everywhere else in the sheet they are inline strings, not
# numbers). Force text format before assigning those specific values so
# they are not auto-coerced into real numbers. Column O is the one
# genuinely-numeric column and is left as a normal numeric assignment.

# --- Row 25 (CBL canonical, TCGA-AB-2914): only the N/O p-value columns change ---
$ws.Range("N25").NumberFormat = "@"
$ws.Range("N25").Value2 = "0.00675675675675676"
$ws.Range("O25").Value2 = 0.006756756756756757

# --- Row 26: becomes the canonical-SJ row for TCGA-AB-2956 sample ---
$ws.Range("A26").Value2 = "TCGA"
$ws.Range("B26").Value2 = "CBL"
$ws.Range("C26").Value2 = "chr11,119278164,A,T"
$ws.Range("D26").Value2 = "ENST00000264033.4:c.1096-2A>T"
$ws.Range("E26").Value2 = "-"
$ws.Range("F26").Value2 = "Aceptor Loss"
$ws.Range("G26").Value2 = "CanonicalSJ"
$ws.Range("H26").Value2 = "chr11:119277845-119278165"
$ws.Range("I26").Value2 = "TCGA-AB-2956"
$ws.Range("J26").Value2 = "TCGA-AB-2956"
$ws.Range("K26").Value2 = "TCGA-AB-2956-03A"
$ws.Range("L26").Value2 = "TCGA-AB-2956-03A"
$ws.Range("M26").NumberFormat = "@"
$ws.Range("M26").Value2 = "2.63157894736842"
$ws.Range("N26").NumberFormat = "@"
$ws.Range("N26").Value2 = "0"
$ws.Range("O26").Value2 = 0
$ws.Range("P26").Value2 = "Not RNA Called"
$ws.Range("Q26").Value2 = "Not DNA Called"
$ws.Range("R26").Value2 = "Not DNA Called"
$ws.Range("S26").Value2 = "Not DNA Called"
$ws.Range("T26").Value2 = "Not DNA Called"

# --- Row 27: Aceptor Gain / AlternativeSJ row for TCGA-AB-2914, shorter SJ ---
$ws.Range("H27").Value2 = "chr11:119277845-119278189"
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value2 = "0.511770726714432"
$ws.Range("N27").NumberFormat = "@"
$ws.Range("N27").Value2 = "1"
$ws.Range("O27").Value2 = 0

# --- Row 28 (new): Aceptor Gain / AlternativeSJ row for TCGA-AB-2956, shorter SJ ---
$ws.Range("A28").Value2 = "TCGA"
$ws.Range("B28").Value2 = "CBL"
$ws.Range("C28").Value2 = "chr11,119278164,A,T"
$ws.Range("D28").Value2 = "ENST00000264033.4:c.1096-2A>T"
$ws.Range("E28").Value2 = "-"
$ws.Range("F28").Value2 = "Aceptor Gain"
$ws.Range("G28").Value2 = "AlternativeSJ found in MUT samples"
$ws.Range("H28").Value2 = "chr11:119277845-119278189"
$ws.Range("I28").Value2 = "TCGA-AB-2956"
$ws.Range("J28").Value2 = "TCGA-AB-2956"
$ws.Range("K28").Value2 = "TCGA-AB-2956-03A"
$ws.Range("L28").Value2 = "TCGA-AB-2956-03A"
$ws.Range("M28").NumberFormat = "@"
$ws.Range("M28").Value2 = "0.657894736842105"
$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value2 = "1"
$ws.Range("O28").Value2 = 0
$ws.Range("P28").Value2 = "Not RNA Called"
$ws.Range("Q28").Value2 = "Not DNA Called"
$ws.Range("R28").Value2 = "Not DNA Called"
$ws.Range("S28").Value2 = "Not DNA Called"
$ws.Range("T28").Value2 = "Not DNA Called"

# --- Row 29 (new): Aceptor Gain / AlternativeSJ row for TCGA-AB-2914, longer SJ ---
$ws.Range("A29").Value2 = "TCGA"
$ws.Range("B29").Value2 = "CBL"
$ws.Range("C29").Value2 = "chr11,119278165,G,C"
$ws.Range("D29").Value2 = "ENST00000264033.4:c.1096-1G>C"
$ws.Range("E29").Value2 = "-"
$ws.Range("F29").Value2 = "Aceptor Gain"
$ws.Range("G29").Value2 = "AlternativeSJ found in MUT samples"
$ws.Range("H29").Value2 = "chr11:119277845-119278237"
$ws.Range("I29").Value2 = "TCGA-AB-2914"
$ws.Range("J29").Value2 = "TCGA-AB-2914"
$ws.Range("K29").Value2 = "TCGA-AB-2914-03A"
$ws.Range("L29").Value2 = "TCGA-AB-2914-03A"
$ws.Range("M29").NumberFormat = "@"
$ws.Range("M29").Value2 = "0.818833162743091"
$ws.Range("N29").NumberFormat = "@"
$ws.Range("N29").Value2 = "1"
$ws.Range("O29").Value2 = 0
$ws.Range("P29").Value2 = "Not RNA Called"
$ws.Range("Q29").NumberFormat = "@"
$ws.Range("Q29").Value2 = "0.278"
$ws.Range("R29").NumberFormat = "@"
$ws.Range("R29").Value2 = "0.3158"
$ws.Range("S29").NumberFormat = "@"
$ws.Range("S29").Value2 = "0.3158"
$ws.Range("T29").NumberFormat = "@"
$ws.Range("T29").Value2 = "0.3158"

# --- Row 30 (new): Aceptor Gain / AlternativeSJ row for TCGA-AB-2956, longer SJ ---
$ws.Range("A30").Value2 = "TCGA"
$ws.Range("B30").Value2 = "CBL"
$ws.Range("C30").Value2 = "chr11,119278164,A,T"
$ws.Range("D30").Value2 = "ENST00000264033.4:c.1096-2A>T"
$ws.Range("E30").Value2 = "-"
$ws.Range("F30").Value2 = "Aceptor Gain"
$ws.Range("G30").Value2 = "AlternativeSJ found in MUT samples"
$ws.Range("H30").Value2 = "chr11:119277845-119278237"
$ws.Range("I30").Value2 = "TCGA-AB-2956"
$ws.Range("J30").Value2 = "TCGA-AB-2956"
$ws.Range("K30").Value2 = "TCGA-AB-2956-03A"
$ws.Range("L30").Value2 = "TCGA-AB-2956-03A"
$ws.Range("M30").NumberFormat = "@"
$ws.Range("M30").Value2 = "0.87719298245614"
$ws.Range("N30").NumberFormat = "@"
$ws.Range("N30").Value2 = "1"
$ws.Range("O30").Value2 = 0
$ws.Range("P30").Value2 = "Not RNA Called"
$ws.Range("Q30").Value2 = "Not DNA Called"
$ws.Range("R30").Value2 = "Not DNA Called"
$ws.Range("S30").Value2 = "Not DNA Called"
$ws.Range("T30").Value2 = "Not DNA Called"

# Rows 31 (PTPN11, ex-row28) and 32 (WT1, ex-row29) already carried their
# original values down via the row insert above, so no further edits are
# needed there.
